$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert 3 new rows before row 176, shifting the existing 176-259 block down to 179-262.
$ws.Rows("176:178").Insert()

# New data for the 3 freshly-inserted rows.
$newRows = @(
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44488, 9, 100112032, "Zapallo italiano", "Bola 8", "Primera", 40, 16000, 16000, 16000, "`$/caja 60 unidades", "Región de O'Higgins", 267, 60, "Hortaliza"),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44488, 9, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 200, 15000, 15000, 15000, "`$/caja 60 unidades", "Región de Arica y Parinacota", 250, 60, "Hortaliza"),
    @(10, "Vega Modelo de Temuco", "La Araucanía", 44488, 9, 100112032, "Zapallo italiano", "Sin especificar", "Primera", 200, 17000, 17000, 17000, "`$/caja 60 unidades", "Región de O'Higgins", 283, 60, "Hortaliza")
)

for ($i = 0; $i -lt 3; $i++) {
    $row = 176 + $i
    $values = $newRows[$i]
    for ($c = 1; $c -le 18; $c++) {
        $ws.Cells.Item($row, $c).Value = $values[$c - 1]
    }
}
